# Applies the job-search tracking updates described by the commit diff:
#  - Adds a custom date number format (yyyy-mm-dd) used for the "date" column
#  - Appends 9 new rows to the "LinkedIn" sheet
#  - Appends 3 new rows to the "Others" sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LinkedIn")
$ws3 = $wb.Worksheets.Item("Others")

$dateSerial = 45092
$dateFormat = "yyyy-mm-dd"

$linkedinRows = @(
    @($dateSerial, "Optum ", "Business Development Consultant (EHR Services) - Remote", "Franklin, Tennessee, United States", "https://www.linkedin.com/jobs/view/business-development-consultant-ehr-services-remote-at-optum-3637581437", "pharmacy ehr manager"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr manager"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr manager"),
    @($dateSerial, "Optum ", "Business Development Consultant (EHR Services) - Remote", "Franklin, Tennessee, United States", "https://www.linkedin.com/jobs/view/business-development-consultant-ehr-services-remote-at-optum-3637581437", "pharmacy ehr specialist"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr specialist"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr specialist"),
    @($dateSerial, "Optum ", "Business Development Consultant (EHR Services) - Remote", "Franklin, Tennessee, United States", "https://www.linkedin.com/jobs/view/business-development-consultant-ehr-services-remote-at-optum-3637581437", "pharmacy ehr associate"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr associate"),
    @($dateSerial, "4 Healing Center ", "Billing Specialist", "Murray, Utah, United States", "https://www.linkedin.com/jobs/view/billing-specialist-at-4-healing-center-3637510897?refId=zUHvkrwS6bcg3flBXyE00w%3D%3D&trackingId=bt71DNAtq%2B7kHr3jn%2Fj5Gw%3D%3D&position=15&pageNum=0&trk=public_jobs_jserp-result_search-card", "pharmacy ehr associate")
)

$othersRows = @(
    @($dateSerial, "none", "none", "none", "https://au.linkedin.com/jobs/view/clinical-coder-admin-off-lvl-3-4-5-or-mra-incremental-fairfield-hospital-casual-at-i-work-for-nsw-3629284842", "pharmacy emr manager"),
    @($dateSerial, "none", "none", "none", "https://au.linkedin.com/jobs/view/clinical-coder-admin-off-lvl-3-4-5-or-mra-incremental-fairfield-hospital-casual-at-i-work-for-nsw-3629284842", "pharmacy emr specialist"),
    @($dateSerial, "none", "none", "none", "https://au.linkedin.com/jobs/view/clinical-coder-admin-off-lvl-3-4-5-or-mra-incremental-fairfield-hospital-casual-at-i-work-for-nsw-3629284842", "pharmacy emr associate")
)

$r = 2
foreach ($row in $linkedinRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$r = 2
foreach ($row in $othersRows) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
